{"js": "const replacements = [\n  [\"2024-03-26 Tuesday\", \"2024-03-27 Wednesday\"],\n  [\"984\u00f74=246, 0\", \"118\u00f75=23, 3\"],\n  [\"298\u00f74=74, 2\", \"769\u00f72=384, 1\"],\n  [\"368\u00f77=52, 4\", \"941\u00f75=188, 1\"],\n  [\"936\u00f74=234, 0\", \"999\u00f78=124, 7\"],\n  [\"830\u00f77=118, 4\", \"251\u00f74=62, 3\"],\n  [\"514\u00f75=102, 4\", \"393\u00f77=56, 1\"],\n  [\"401\u00f79=44, 5\", \"977\u00f72=488, 1\"],\n  [\"320\u00f77=45, 5\", \"962\u00f77=137, 3\"],\n  [\"964\u00f73=321, 1\", \"203\u00f75=40, 3\"],\n  [\"867\u00f78=108, 3\", \"994\u00f76=165, 4\"],\n  [\"518\u00f72=259, 0\", \"218\u00f72=109, 0\"],\n  [\"334\u00f74=83, 2\", \"890\u00f72=445, 0\"],\n  [\"931\u00f75=186, 1\", \"584\u00f74=146, 0\"],\n  [\"618\u00f72=309, 0\", \"275\u00f79=30, 5\"],\n  [\"589\u00f76=98, 1\", \"728\u00f74=182, 0\"],\n  [\"356\u00f76=59, 2\", \"280\u00f75=56, 0\"],\n  [\"810\u00f78=101, 2\", \"940\u00f78=117, 4\"],\n  [\"943\u00f72=471, 1\", \"722\u00f74=180, 2\"],\n  [\"992\u00f74=248, 0\", \"820\u00f72=410, 0\"],\n  [\"520\u00f79=57, 7\", \"262\u00f77=37, 3\"],\n  [\"164\u00f73=54, 2\", \"626\u00f78=78, 2\"],\n  [\"112\u00f79=12, 4\", \"741\u00f74=185, 1\"],\n  [\"115\u00f78=14, 3\", \"681\u00f77=97, 2\"],\n  [\"719\u00f72=359, 1\", \"754\u00f75=150, 4\"],\n  [\"550\u00f77=78, 4\", \"439\u00f72=219, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-03-26 Tuesday\"; New = \"2024-03-27 Wednesday\" },\n    @{ Old = \"984\u00f74=246, 0\"; New = \"118\u00f75=23, 3\" },\n    @{ Old = \"298\u00f74=74, 2\"; New = \"769\u00f72=384, 1\" },\n    @{ Old = \"368\u00f77=52, 4\"; New = \"941\u00f75=188, 1\" },\n    @{ Old = \"936\u00f74=234, 0\"; New = \"999\u00f78=124, 7\" },\n    @{ Old = \"830\u00f77=118, 4\"; New = \"251\u00f74=62, 3\" },\n    @{ Old = \"514\u00f75=102, 4\"; New = \"393\u00f77=56, 1\" },\n    @{ Old = \"401\u00f79=44, 5\"; New = \"977\u00f72=488, 1\" },\n    @{ Old = \"320\u00f77=45, 5\"; New = \"962\u00f77=137, 3\" },\n    @{ Old = \"964\u00f73=321, 1\"; New = \"203\u00f75=40, 3\" },\n    @{ Old = \"867\u00f78=108, 3\"; New = \"994\u00f76=165, 4\" },\n    @{ Old = \"518\u00f72=259, 0\"; New = \"218\u00f72=109, 0\" },\n    @{ Old = \"334\u00f74=83, 2\"; New = \"890\u00f72=445, 0\" },\n    @{ Old = \"931\u00f75=186, 1\"; New = \"584\u00f74=146, 0\" },\n    @{ Old = \"618\u00f72=309, 0\"; New = \"275\u00f79=30, 5\" },\n    @{ Old = \"589\u00f76=98, 1\"; New = \"728\u00f74=182, 0\" },\n    @{ Old = \"356\u00f76=59, 2\"; New = \"280\u00f75=56, 0\" },\n    @{ Old = \"810\u00f78=101, 2\"; New = \"940\u00f78=117, 4\" },\n    @{ Old = \"943\u00f72=471, 1\"; New = \"722\u00f74=180, 2\" },\n    @{ Old = \"992\u00f74=248, 0\"; New = \"820\u00f72=410, 0\" },\n    @{ Old = \"520\u00f79=57, 7\"; New = \"262\u00f77=37, 3\" },\n    @{ Old = \"164\u00f73=54, 2\"; New = \"626\u00f78=78, 2\" },\n    @{ Old = \"112\u00f79=12, 4\"; New = \"741\u00f74=185, 1\" },\n    @{ Old = \"115\u00f78=14, 3\"; New = \"681\u00f77=97, 2\" },\n    @{ Old = \"719\u00f72=359, 1\"; New = \"754\u00f75=150, 4\" },\n    @{ Old = \"550\u00f77=78, 4\"; New = \"439\u00f72=219, 1\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
